# Naive baseline that analyzes measures independently:
# halve the step size of the Properties column (A2:A26),
# going from steps of 0.1 to steps of 0.05.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 2; $i -le 26; $i++) {
    $cell = $ws.Cells.Item($i, 1)
    $cell.Value = $cell.Value() / 2
}

$ws.Range("B19").Select()
